# ----------------------------------------------------------------------------
# Edit: rewrite the body of the report.
#
# The single "In order to handle the input string..." paragraph is split into
# a new "Overview:" heading + two descriptive paragraphs about the project,
# a blank spacer paragraph, a new bold "Handling the input string:" heading
# (carrying the existing _GoBack bookmark + a lastRenderedPageBreak marker),
# and the original paragraph (still starting "In order to handle..." but now
# with a "from" grammar-checker proofErr pair inserted and its firstLine
# indent retained).
#
# Because this reshuffles/re-splits runs, proofErr markers, a tab character,
# bookmarks and a page-break marker across many new paragraphs, it is far more
# reliable to express the whole new body as one WordprocessingML fragment and
# hand it to Word via Range.InsertXML (the Flat OPC "single XML file" format
# Word itself uses for clipboard / InsertXML payloads) rather than to
# choreograph dozens of individual Find/Replace & InsertParagraph calls.
# ----------------------------------------------------------------------------

$d = $word.ActiveDocument

# The replacement WordprocessingML for every paragraph in the body (everything
# except the trailing <w:sectPr>, which Word keeps from the existing document).
$newBodyParagraphs = @'
<w:p w:rsidR="00404501" w:rsidRPr="00E77260" w:rsidRDefault="00B551A8" w:rsidP="00E77260"><w:pPr><w:spacing w:line="600" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00E77260"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Report</w:t></w:r></w:p><w:p w:rsidR="00B551A8" w:rsidRPr="00E77260" w:rsidRDefault="00B551A8" w:rsidP="00E77260"><w:pPr><w:spacing w:line="600" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00E77260"><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>Overview:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="600" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/><w:t xml:space="preserve">The project is a </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">CMOS circuit spice </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>netlist code generator</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>. The user will input a Boolean expression made of and, or, and not in addition to normal alphabets. The output is supposed to be the netlist code for this circuit consisting of any CMOS transistor (NMOS and PMOS). The circuit is divided into two parts: the pull up circuit (PUN) and the pull down circuit (PDN), with the PUN made of PMOS transistor, while the PDN constituted of NMOS transistors. Our implementation language is C++ and each line of code is a transistor with all its junctions which are its name, drain, source, gate, body, and finally its type.</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="600" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:tab/><w:t xml:space="preserve">The bonus feature we implemented was the first one which generating a full spice deck compatible with the free </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>LTSpice</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">. Thus, the code generated is not just the data statements; instead it is generates a title statement, end statement, control statements, and output statements in addition to the data statements of course. </w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="600" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:spacing w:line="600" w:lineRule="auto"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:lastRenderedPageBreak/><w:t>Handling the input string</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>:</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="600" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">In order to handle the input string efficiently, we thought we could change it </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>from</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> its infix form to prefix. For this purpose, stacks were used where we pushed and popped parts of the strings depending on the type of the character. Also, in handling the string, we were aware of the fact that pull up circuit should in terms of complemented input and pull down are the opposite. Therefore, </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">for the pull up for example, </w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>i</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>f an input was not complemented, we complemented and if it was complemented we left as it is.</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> Finally, the output prefix was a single string so we had to divide the string into an array of strings so that we could deal with every element in the array on i</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>ts</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> own</w:t></w:r><w:r><w:rPr><w:rFonts w:asciiTheme="majorBidi" w:hAnsiTheme="majorBidi" w:cstheme="majorBidi"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r></w:p>
'@

$flatOpc = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><?xml version="1.0" encoding="UTF-8" standalone="yes"?><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14"><w:body>__PLACEHOLDER__<w:sectPr><w:pgSz w:w="12240" w:h="15840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$flatOpc = $flatOpc.Replace('__PLACEHOLDER__', $newBodyParagraphs)

# Replace the entire document body content in one shot (equivalent to
# selecting Ctrl+A and pasting the new OOXML over it); Word keeps the
# document's own final section properties (sectPr) regardless of what the
# inserted fragment's own sectPr says.
$d.Content.InsertXML($flatOpc, "Replace")
